# Update TPM-derived NATMI metrics for Vegfa-Kdr LR-pair sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.436534333333333
$ws.Range("H2").Value = 4.309603
$ws.Range("I2").Value = 0.03241561610838976
$ws.Range("J2").Value = 0.03241561610838976
$ws.Range("M2").Value = 174.1282373333333
$ws.Range("N2").Value = 522.384712
$ws.Range("O2").Value = 0.985625830323027
$ws.Range("P2").Value = 0.985625830323027
$ws.Range("Q2").Value = 250.1411913321485
$ws.Range("R2").Value = 2251.270721989336
$ws.Range("S2").Value = 0.03194966854226414
$ws.Range("T2").Value = 0.03194966854226414

# Row 3
$ws.Range("G3").Value = 1.436534333333333
$ws.Range("H3").Value = 4.309603
$ws.Range("I3").Value = 0.03241561610838976
$ws.Range("J3").Value = 0.03241561610838976
$ws.Range("O3").Value = 0.003686901313133159
$ws.Range("P3").Value = 0.003686901313133159
$ws.Range("Q3").Value = 0.9356957360674444
$ws.Range("R3").Value = 8.421261624607
$ws.Range("S3").Value = 0.0001195131775960426
$ws.Range("T3").Value = 0.0001195131775960426

# Row 4
$ws.Range("G4").Value = 1.436534333333333
$ws.Range("H4").Value = 4.309603
$ws.Range("I4").Value = 0.03241561610838976
$ws.Range("J4").Value = 0.03241561610838976
$ws.Range("M4").Value = 1.888095
$ws.Range("N4").Value = 5.664285
$ws.Range("O4").Value = 0.01068726836383999
$ws.Range("P4").Value = 0.01068726836383999
$ws.Range("Q4").Value = 2.712313292095
$ws.Range("R4").Value = 24.410819628855
$ws.Range("S4").Value = 0.000346434388529576
$ws.Range("T4").Value = 0.000346434388529576

# Row 5
$ws.Range("G5").Value = 34.88211266666666
$ws.Range("I5").Value = 0.7871201871162607
$ws.Range("J5").Value = 0.7871201871162609
$ws.Range("M5").Value = 174.1282373333333
$ws.Range("N5").Value = 522.384712
$ws.Range("O5").Value = 0.985625830323027
$ws.Range("P5").Value = 0.985625830323027
$ws.Range("Q5").Value = 6073.960793109406
$ws.Range("R5").Value = 54665.64713798466
$ws.Range("S5").Value = 0.7758059879904808
$ws.Range("T5").Value = 0.7758059879904809

# Row 6
$ws.Range("G6").Value = 34.88211266666666
$ws.Range("I6").Value = 0.7871201871162607
$ws.Range("J6").Value = 0.7871201871162609
$ws.Range("O6").Value = 0.003686901313133159
$ws.Range("P6").Value = 0.003686901313133159
$ws.Range("S6").Value = 0.002902034451472559
$ws.Range("T6").Value = 0.00290203445147256

# Row 7
$ws.Range("G7").Value = 34.88211266666666
$ws.Range("I7").Value = 0.7871201871162607
$ws.Range("J7").Value = 0.7871201871162609
$ws.Range("M7").Value = 1.888095
$ws.Range("N7").Value = 5.664285
$ws.Range("O7").Value = 0.01068726836383999
$ws.Range("P7").Value = 0.01068726836383999
$ws.Range("Q7").Value = 65.86074251536999
$ws.Range("R7").Value = 592.7466826383301
$ws.Range("S7").Value = 0.008412164674307431
$ws.Range("T7").Value = 0.008412164674307431

# Row 8
$ws.Range("G8").Value = 7.997472999999999
$ws.Range("H8").Value = 23.992419
$ws.Range("I8").Value = 0.1804641967753495
$ws.Range("J8").Value = 0.1804641967753495
$ws.Range("M8").Value = 174.1282373333333
$ws.Range("N8").Value = 522.384712
$ws.Range("O8").Value = 0.985625830323027
$ws.Range("P8").Value = 0.985625830323027
$ws.Range("Q8").Value = 1392.585876610925
$ws.Range("R8").Value = 12533.27288949833
$ws.Range("S8").Value = 0.1778701737902819
$ws.Range("T8").Value = 0.1778701737902819

# Row 9
$ws.Range("G9").Value = 7.997472999999999
$ws.Range("H9").Value = 23.992419
$ws.Range("I9").Value = 0.1804641967753495
$ws.Range("J9").Value = 0.1804641967753495
$ws.Range("O9").Value = 0.003686901313133159
$ws.Range("P9").Value = 0.003686901313133159
$ws.Range("Q9").Value = 5.209204689212333
$ws.Range("R9").Value = 46.882842202911
$ws.Range("S9").Value = 0.0006653536840645567
$ws.Range("T9").Value = 0.0006653536840645567

# Row 10
$ws.Range("G10").Value = 7.997472999999999
$ws.Range("H10").Value = 23.992419
$ws.Range("I10").Value = 0.1804641967753495
$ws.Range("J10").Value = 0.1804641967753495
$ws.Range("M10").Value = 1.888095
$ws.Range("N10").Value = 5.664285
$ws.Range("O10").Value = 0.01068726836383999
$ws.Range("P10").Value = 0.01068726836383999
$ws.Range("Q10").Value = 15.099988783935
$ws.Range("R10").Value = 135.899899055415
$ws.Range("S10").Value = 0.001928669301002988
$ws.Range("T10").Value = 0.001928669301002988
